# Assignment5.xlsx edit script
# 1) "Measures" sheet / Table1: insert a new "Measure Folder" column before
#    "Measure Description" (so the table grows from A1:D4 to A1:E4), fill the
#    new column with "No Folder Defined", reformat the DAX expressions in
#    column B, and rewrite the descriptions in the (now) E column.
# 2) "Source Information" sheet / Table2: bump the "Table No" for the existing
#    LiteracyTable row from 1 to 2, and reword the Modification Description.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Measures
$ws2 = $wb.Worksheets.Item(2)   # Source Information

# ---------------------------------------------------------------------
# Measures sheet: shift column D ("Measure Description") data into E,
# keep formatting consistent with the existing columns.
# ---------------------------------------------------------------------

# Copy column D formatting into the about-to-be-used column E (single cells,
# not whole-column, so we don't blow the used range out to 1,048,576 rows).
$ws1.Range("D1").Copy()
$ws1.Range("E1").PasteSpecial(-4122)
$ws1.Range("D2").Copy()
$ws1.Range("E2").PasteSpecial(-4122)
$ws1.Range("D3").Copy()
$ws1.Range("E3").PasteSpecial(-4122)
$ws1.Range("D4").Copy()
$ws1.Range("E4").PasteSpecial(-4122)

# Move "Measure Description" text into column E.
$ws1.Range("E1").Value = "Measure Description"
$ws1.Range("E2").Value = "The DISTINCTCOUNT function in LiteracyTable[State] is used to count the unique occurrences of states in the table. This helps provide an understanding of how many different states are represented in the Literacy Table."
$ws1.Range("E3").Value = "This calculation counts the number of cities in the LiteracyTable that have a Level value of ""City""."
$ws1.Range("E4").Value = "This calculation counts the number of unique cities in a literacy table for all records with a Level of ""UA"". It tells us how many cities have a literacy level of ""UA""."

# Column D now becomes "Measure Folder".
$ws1.Range("D1").Value = "Measure Folder"
$ws1.Range("D2").Value = "No Folder Defined"
$ws1.Range("D3").Value = "No Folder Defined"
$ws1.Range("D4").Value = "No Folder Defined"

# Reformat the DAX expressions in column B (pretty-printed, multi-line).
$ws1.Range("B2").Value = "`nDISTINCTCOUNT(LiteracyTable[State])"
$ws1.Range("B3").Value = "`nCOUNTAX(`n    FILTER(`n        LiteracyTable,`n        'LiteracyTable'[Level] = ""City""`n    ), `n    [City]`n)"
$ws1.Range("B4").Value = "`nCOUNTAX(`n    FILTER(`n        LiteracyTable, `n        LiteracyTable[Level] = ""UA""`n    ), `n    [City]`n)"

# Resize the ListObject/table so it spans the new column and rename the
# header cells through the table's header range (keeps table1.xml in sync).
$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:E4"))
$lo1.HeaderRowRange.Cells.Item(1, 4).Value = "Measure Folder"
$lo1.HeaderRowRange.Cells.Item(1, 5).Value = "Measure Description"

# Column widths: D keeps 30 (was 50), new column E gets 50.
$ws1.Columns.Item(4).ColumnWidth = 29.1
$ws1.Columns.Item(5).ColumnWidth = 49.1

# ---------------------------------------------------------------------
# Source Information sheet: LiteracyTable's "Table No" goes from 1 to 2,
# and the Modification Description (I2) is reworded.
# ---------------------------------------------------------------------

$ws2.Range("A2").Value = 2

$ws2.Range("I2").Value = "1. This changes the data types of the columns in the LiteracyTable_Table to integers, text, or numbers.`n2. This sentence means that the ""Replaced Errors"" table was created by changing the type and replacing the ""State Code"" value with 33.`n3. Renamed Columns means that the column labelled ""Name of Urban Agglomeration/City"" has been changed to ""City"" for the table ""Replaced Errors"".`n4. Split a ""City"" column into two columns named ""City.1"" and ""City.2"" based on spaces between words.`n5. Table.TransformColumnTypes takes a data table as an input and transforms the type of certain columns from one type to another.`n6. This line of code removes the ""City.2"" column from the table ""Changed Type1"".`n7. The column names ""City.1"" and ""State Name"" were changed to ""City"" and ""State"", respectively.`n"
